$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = 99
$ws.Range("B12").Value = "ALTRO MOTIVO DI ESENZIONE"

$ws.Range("A12:B12").VerticalAlignment = -4108
$ws.Range("A12:B12").WrapText = $true

$ws.Activate()
$ws.Range("C12").Select()

